$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily records (data aggiornati fino al 20/09/2021)
$data = @(
    @(375, 44449, 6, 24, 59.58587814687919),
    @(376, 44450, 0, 20, 49.65489845573266),
    @(377, 44451, 4, 23, 57.10313322409256),
    @(378, 44452, 0, 19, 47.17215353294603),
    @(379, 44453, 4, 18, 44.68940861015939),
    @(380, 44454, 0, 16, 39.72391876458613),
    @(381, 44455, 9, 23, 57.10313322409256),
    @(382, 44456, 2, 19, 47.17215353294603),
    @(383, 44457, 2, 21, 52.13764337851929),
    @(384, 44458, 1, 18, 44.68940861015939),
    @(385, 44459, 2, 20, 49.65489845573266)
)

$lastRow = 374

foreach ($row in $data) {
    $r = $row[0]

    # Copy the date cell above so the new row inherits the same style (s="2",
    # bordered/centered/bold with the YYYY-MM-DD HH:MM:SS number format)
    # instead of picking up a brand-new auto-generated style.
    $ws.Range("A$lastRow").Copy($ws.Range("A$r"))

    $ws.Range("A$r").Value = $row[1]
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
    $ws.Range("D$r").Value = $row[4]

    $lastRow = $r
}
